# Update OPL with tasks from Kick-Off-Meeting
$wb = $excel.ActiveWorkbook

$wsOpl = $wb.Worksheets.Item("OPL Vorlage")
$wsTab = $wb.Worksheets.Item("Tabelle1")

# --- OPL Vorlage: fill in rows 4 and 5 with new open points from the Kick-Off-Meeting ---
$wsOpl.Range("B4").Value = "POC der APIs"
$wsOpl.Range("C4").Value = "Marius Kurth"
$wsOpl.Range("E4").Value = "Details s. Protokoll 15.05.2024"
$wsOpl.Range("D4").Value = "bis 24.05.2024"
$wsOpl.Range("G4").Value = "A"
$wsOpl.Range("H4").Value = "iB"

$wsOpl.Range("B5").Value = "Design-Entwurf"
$wsOpl.Range("C5").Value = "Matteo Kosina"
$wsOpl.Range("D5").Value = "bis 24.05.2024"
$wsOpl.Range("G5").Value = "W"
$wsOpl.Range("H5").Value = "iB"

# --- Tabelle1: replace the "Verantwortlich" lookup list with the meeting's participants ---
$wsTab.Range("C3").Value = "Leon Fertig"
$wsTab.Range("C4").Value = "Matteo Kosina"
$wsTab.Range("C5").Value = "Marius Kurth"
$wsTab.Range("C6").Value = ""
$wsTab.Range("C7").Value = ""

# --- Update the saved selections/active cells to match the author's final cursor position ---
$wsTab.Range("C12").Select() | Out-Null
$wsOpl.Activate() | Out-Null
$wsOpl.Range("B12").Select() | Out-Null
